# Add a new "2022-Q1" sheet (fund holdings detail, same shape as the
# other quarterly sheets) right before the "总计" (totals) summary sheet,
# and prepend a corresponding "2022-Q1" row to the "总计" sheet's table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet just before "总计".
#    NOTE: sheet variables here are positional, not stable object
#    handles -- after Worksheets.Add() shifts sheets around, previously
#    captured references can now point at a different sheet. So we
#    re-fetch "总计" by name whenever we need it after a mutation.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# The header row (B1:H1) and the index column (A2:A5) use the same
# bold/bordered style as the equivalent cells on the other quarterly
# sheets -- copy it across since a brand-new sheet starts with no such
# style defined.
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows. Fund codes / size / position figures are stored as TEXT in
# the source workbook (not numbers) -- force the number format to text
# first so leading zeros (e.g. "006679") and the literal string
# representation survive instead of being coerced to a number. The
# number format is stripped back off afterwards (ClearFormats keeps the
# stored text value/type, it just drops the style record) so these
# cells end up as plain unstyled text cells like the source sheet.
$newSheet.Range("B2:G5").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Range("B2").Value = "006679"
$newSheet.Range("C2").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A"
$newSheet.Range("D2").Value = "14.75"
$newSheet.Range("E2").Value = "83.19"
$newSheet.Range("F2").Value = "7.95"
$newSheet.Range("G2").Value = "1.1726"
$newSheet.Range("H2").Value = 2

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Range("B3").Value = "162719"
$newSheet.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
$newSheet.Range("D3").Value = "14.75"
$newSheet.Range("E3").Value = "83.19"
$newSheet.Range("F3").Value = "7.95"
$newSheet.Range("G3").Value = "1.1726"
$newSheet.Range("H3").Value = 2

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Range("B4").Value = "006680"
$newSheet.Range("C4").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C"
$newSheet.Range("D4").Value = "4.73"
$newSheet.Range("E4").Value = "83.19"
$newSheet.Range("F4").Value = "7.95"
$newSheet.Range("G4").Value = "0.3760"
$newSheet.Range("H4").Value = 2

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Range("B5").Value = "004243"
$newSheet.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
$newSheet.Range("D5").Value = "4.73"
$newSheet.Range("E5").Value = "83.19"
$newSheet.Range("F5").Value = "7.95"
$newSheet.Range("G5").Value = "0.3760"
$newSheet.Range("H5").Value = 2

# Drop the temporary text number-format/style now that the values are
# safely stored as strings (leaves these cells with no style, matching
# the other quarterly sheets).
$newSheet.Range("B2:G5").ClearFormats()

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" sheet, pushing existing rows
#    down by one and renumbering the leading index column.
#    Re-fetch "总计" by name since the earlier Add() shifted positions.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Insert() clones row 1's (header) formatting into the new row 2 cells;
# strip that back off so B2:D2 end up unstyled like the other data rows.
$totalSheet.Range("B2:D2").ClearFormats()

# A2 needs the same style as the other index-column cells (A3:A7) --
# copy it across explicitly since Insert() left A2 completely unstyled.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 3.1

# Renumber the remaining index column (A) 1..5 for the rows that shifted
# down (previously 0..4 for 2021-Q4 .. 2020-Q4).
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
